$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Save" header in H1, matching the formatting of the other header cells (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell H2
$ws.Range("H2").Value = 0
